# recup du form et injection id nom prenom ok
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3
$ws.Range("B1").Value = "Cherief"
$ws.Range("C1").Value = "Saufiane"

$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "toto"
$ws.Range("C2").Value = "tata"
